$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared "Colorcode" values:
#   #445566 -> #f54266 (rows with "good" storability)
#   #993422 -> #42f548 (rows with "bad" storability)
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -eq "#445566") {
        $cell.Value = "#f54266"
    } elseif ($val -eq "#993422") {
        $cell.Value = "#42f548"
    }
}

# Move the active selection from D10 to B2
$ws.Range("B2").Select()
